$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.4627908493242501
$ws.Range("C2").Value = 0.1134005571835246
$ws.Range("D2").Value = 0.07843944507473566
$ws.Range("E2").Value = 0.41152745995565
$ws.Range("G2").Value = 0.9516726143155552
$ws.Range("H2").Value = 0.9464481696558238
$ws.Range("I2").Value = 0.7481229575361183
$ws.Range("K2").Value = 0.5578637603159393
$ws.Range("N2").Value = 1.605780702952011

$ws.Range("B3").Value = 0.4164782384311252
$ws.Range("C3").Value = 0.100008059534531
$ws.Range("D3").Value = 0.07115579397002136
$ws.Range("E3").Value = 0.3589963012006763
$ws.Range("G3").Value = 0.9357387779456019
$ws.Range("H3").Value = 0.9442726043521219
$ws.Range("I3").Value = 0.7457616126376934
$ws.Range("K3").Value = 0.4991180516498446
$ws.Range("N3").Value = 1.620824802326169

$ws.Range("B4").Value = 0.3882077155305126
$ws.Range("C4").Value = 0.09181565796373548
$ws.Range("D4").Value = 0.06672095575694925
$ws.Range("E4").Value = 0.3268459619217765
$ws.Range("G4").Value = 0.9266061884878667
$ws.Range("H4").Value = 0.9434249534876358
$ws.Range("I4").Value = 0.7447370193433471
$ws.Range("K4").Value = 0.4632373810664774
$ws.Range("N4").Value = 1.630640492642875

$ws.Range("B5").Value = 0.3767287265012271
$ws.Range("C5").Value = 0.08848466053643733
$ws.Range("D5").Value = 0.06492304570784313
$ws.Range("E5").Value = 0.3137685166095849
$ws.Range("G5").Value = 0.9230474580757146
$ws.Range("H5").Value = 0.9432020387869926
$ws.Range("I5").Value = 0.7444261187025702
$ws.Range("K5").Value = 0.4486629244394749
$ws.Range("N5").Value = 1.634785652224053

$ws.Range("B6").Value = 0.3748251497044919
$ws.Range("C6").Value = 0.08793199551664088
$ws.Range("D6").Value = 0.064625065143602
$ws.Range("E6").Value = 0.3115984107982968
$ws.Range("G6").Value = 0.9224663478849067
$ws.Range("H6").Value = 0.9431724163097357
$ws.Range("I6").Value = 0.7443809250063893
$ws.Range("K6").Value = 0.4462456860477744
$ws.Range("N6").Value = 1.63548271311231

$ws.Range("B7").Value = 0.3880527381378442
$ws.Range("C7").Value = 0.09177070501837647
$ws.Range("D7").Value = 0.0666966708585619
$ws.Range("E7").Value = 0.3266695001385642
$ws.Range("G7").Value = 0.926557535781825
$ws.Range("H7").Value = 0.9434214514703427
$ws.Range("I7").Value = 0.7447323951074836
$ws.Range("K7").Value = 0.4630406342605227
$ws.Range("N7").Value = 1.630695808266161

$ws.Range("B8").Value = 0.4467877734882109
$ws.Range("C8").Value = 0.1087763182613912
$ws.Range("D8").Value = 0.07592024481077431
$ws.Range("E8").Value = 0.3933918631780955
$ws.Range("G8").Value = 0.9460430603804468
$ws.Range("H8").Value = 0.9455965692559829
$ws.Range("I8").Value = 0.747220307773226
$ws.Range("K8").Value = 0.5375686460493512
$ws.Range("N8").Value = 1.610847652416552

$ws.Range("B9").Value = 0.5632948942680684
$ws.Range("C9").Value = 0.1423796837788416
$ws.Range("D9").Value = 0.09430818116064188
$ws.Range("E9").Value = 0.5251593183205188
$ws.Range("G9").Value = 0.9894572401785808
$ws.Range("H9").Value = 0.953747636528135
$ws.Range("I9").Value = 0.7554888540765319
$ws.Range("K9").Value = 0.6852481169212012
$ws.Range("N9").Value = 1.576528957162708

$ws.Range("B10").Value = 0.6497307048636003
$ws.Range("C10").Value = 0.167243209556176
$ws.Range("D10").Value = 0.1080081951953105
$ws.Range("E10").Value = 0.6226802770422921
$ws.Range("G10").Value = 1.024584378090708
$ws.Range("H10").Value = 0.9621243742509193
$ws.Range("I10").Value = 0.763653352927534
$ws.Range("K10").Value = 0.7947303019773244
$ws.Range("N10").Value = 1.554138721778259

$ws.Range("B11").Value = 0.689241366420049
$ws.Range("C11").Value = 0.178596672271766
$ws.Range("D11").Value = 0.1142836718732241
$ws.Range("E11").Value = 0.6672322013651524
$ws.Range("G11").Value = 1.041279229859441
$ws.Range("H11").Value = 0.9664581177512446
$ws.Range("I11").Value = 0.7678264921839428
$ws.Range("K11").Value = 0.8447612546671621
$ws.Range("N11").Value = 1.544569764849548

$ws.Range("B12").Value = 0.7042307926226385
$ws.Range("C12").Value = 0.1829024107094028
$ws.Range("D12").Value = 0.1166663569413515
$ws.Range("E12").Value = 0.6841325527182818
$ws.Range("G12").Value = 1.047704980494188
$ws.Range("H12").Value = 0.9681747447642692
$ws.Range("I12").Value = 0.7694731578590535
$ws.Range("K12").Value = 0.8637400029217019
$ws.Range("N12").Value = 1.541035245077964

$ws.Range("B13").Value = 0.7010013210150419
$ws.Range("C13").Value = 0.1819748028632091
$ws.Range("D13").Value = 0.1161529217471013
$ws.Range("E13").Value = 0.6804914104346693
$ws.Range("G13").Value = 1.046316450251425
$ws.Range("H13").Value = 0.9678016741358704
$ws.Range("I13").Value = 0.769115560277406
$ws.Range("K13").Value = 0.8596511095330186
$ws.Range("N13").Value = 1.541792500973699

$ws.Range("B14").Value = 0.6904740013417552
$ws.Range("C14").Value = 0.1789507776440757
$ws.Range("D14").Value = 0.1144795701434447
$ws.Range("E14").Value = 0.6686219988179971
$ws.Range("G14").Value = 1.041805795241231
$ws.Range("H14").Value = 0.966597830124897
$ws.Range("I14").Value = 0.7679606316731551
$ws.Range("K14").Value = 0.8463219813214948
$ws.Range("N14").Value = 1.544277191110538

$ws.Range("B15").Value = 0.6840293137224194
$ws.Range("C15").Value = 0.1770993197407051
$ws.Range("D15").Value = 0.1134554164548689
$ws.Range("E15").Value = 0.661355552756703
$ws.Range("G15").Value = 1.039056430160514
$ws.Range("H15").Value = 0.9658702866840656
$ws.Range("I15").Value = 0.7672618614942337
$ws.Range("K15").Value = 0.8381618319402548
$ws.Range("N15").Value = 1.545810742378066

$ws.Range("B16").Value = 0.6471524412946792
$ws.Range("C16").Value = 0.1665021195885288
$ws.Range("D16").Value = 0.1075989552324756
$ws.Range("E16").Value = 0.6197727065653709
$ws.Range("G16").Value = 1.023507802097811
$ws.Range("H16").Value = 0.9618517066989511
$ws.Range("I16").Value = 0.7633898980317184
$ws.Range("K16").Value = 0.7914652689639752
$ws.Range("N16").Value = 1.554776513703104

$ws.Range("B17").Value = 0.6245786196502934
$ws.Range("C17").Value = 0.160012254913255
$ws.Range("D17").Value = 0.1040173425235906
$ws.Range("E17").Value = 0.5943130540709234
$ws.Range("G17").Value = 1.014153080884171
$ws.Range("H17").Value = 0.9595206398803668
$ws.Range("I17").Value = 0.7611324174905576
$ws.Range("K17").Value = 0.7628769045652746
$ws.Range("N17").Value = 1.560434892078078

$ws.Range("B18").Value = 0.611612680646715
$ws.Range("C18").Value = 0.156283480687847
$ws.Range("D18").Value = 0.1019613619572084
$ws.Range("E18").Value = 0.5796869995335499
$ws.Range("G18").Value = 1.008839769069368
$ws.Range("H18").Value = 0.958229092039204
$ws.Range("I18").Value = 0.7598771613311683
$ws.Range("K18").Value = 0.746454966667784
$ws.Range("N18").Value = 1.563747443482967

$ws.Range("B19").Value = 0.6072257096211899
$ws.Range("C19").Value = 0.1550216652181291
$ws.Range("D19").Value = 0.1012659387095738
$ws.Range("E19").Value = 0.574737832883855
$ws.Range("G19").Value = 1.00705230027134
$ws.Range("H19").Value = 0.9578002404364838
$ws.Range("I19").Value = 0.7594595592048705
$ws.Range("K19").Value = 0.7408984303346529
$ws.Range("N19").Value = 1.564878966592502

$ws.Range("B20").Value = 0.6269797848495671
$ws.Range("C20").Value = 0.1607026937369938
$ws.Range("D20").Value = 0.1043981897046962
$ws.Range("E20").Value = 0.5970214327220305
$ws.Range("G20").Value = 1.015141937552102
$ws.Range("H20").Value = 0.9597636896341442
$ws.Range("I20").Value = 0.7613682577786776
$ws.Range("K20").Value = 0.7659179734388317
$ws.Range("N20").Value = 1.559826542309828

$ws.Range("B21").Value = 0.6935653797994519
$ws.Range("C21").Value = 0.1798388306231971
$ws.Range("D21").Value = 0.1149709025310557
$ws.Range("E21").Value = 0.6721075134695127
$ws.Range("G21").Value = 1.043127860317668
$ws.Range("H21").Value = 0.966949375883928
$ws.Range("I21").Value = 0.7682980574254259
$ws.Range("K21").Value = 0.8502361668726053
$ws.Range("N21").Value = 1.543544957834776

$ws.Range("B22").Value = 0.7372438613971326
$ws.Range("C22").Value = 0.1923829942354587
$ws.Range("D22").Value = 0.1219175284262519
$ws.Range("E22").Value = 0.7213536737698121
$ws.Range("G22").Value = 1.062023572976045
$ws.Range("H22").Value = 0.9720860446286679
$ws.Range("I22").Value = 0.7732142219275389
$ws.Range("K22").Value = 0.9055362132981486
$ws.Range("N22").Value = 1.533423158863677

$ws.Range("B23").Value = 0.7139170418475373
$ws.Range("C23").Value = 0.1856844126381816
$ws.Range("D23").Value = 0.1182065957978722
$ws.Range("E23").Value = 0.6950534662712045
$ws.Range("G23").Value = 1.051882888634964
$ws.Range("H23").Value = 0.9693041074282007
$ws.Range("I23").Value = 0.7705548239716151
$ws.Range("K23").Value = 0.876003696719124
$ws.Range("N23").Value = 1.538777720616764

$ws.Range("B24").Value = 0.6258941805396887
$ws.Range("C24").Value = 0.1603905390917362
$ws.Range("D24").Value = 0.1042259989317387
$ws.Range("E24").Value = 0.5957969405463643
$ws.Range("G24").Value = 1.014694673542749
$ws.Range("H24").Value = 0.9596536554321915
$ws.Range("I24").Value = 0.7612615017164472
$ws.Range("K24").Value = 0.7645430629644636
$ws.Range("N24").Value = 1.560101392003986

$ws.Range("B25").Value = 0.5316315519309853
$ws.Range("C25").Value = 0.1332597433565752
$ws.Range("D25").Value = 0.08930084358107138
$ws.Range("E25").Value = 0.489398547564889
$ws.Range("G25").Value = 0.9771494820024458
$ws.Range("H25").Value = 0.9511245888112967
$ws.Range("I25").Value = 0.7528867236163492
$ws.Range("K25").Value = 0.6451280359810028
$ws.Range("N25").Value = 1.5853184725703
